$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every timestamp in column A (rows 2-97) forward by 2 days
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 2
}

# New "Actual Production (MW)" values for rows 2-31 (rows 32-97 remain 0/unchanged)
$bValues = @(194,188,181,170,156,148,143,140,139,136,137,127,121,114,111,107,111,107,102,98,91,86,84,77,72,70,71,73,79,0)

for ($i = 0; $i -lt $bValues.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 2).Value2 = $bValues[$i]
}
